$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("1. General")
$ws2 = $wb.Worksheets.Item("2. SPS - Area Specification")
$ws3 = $wb.Worksheets.Item("3. SPS - Generic Parameters")

# --- Sheet 3 "3. SPS - Generic Parameters": fix truncated-to-integer values ---
# Row 2
$ws3.Range("B2").Value = 0.1
$ws3.Range("C2").Value = 0
$ws3.Range("D2").Value = 66.7
$ws3.Range("F2").Value = 2
$ws3.Range("G2").Value = 2
$ws3.Range("H2").Value = 1
$ws3.Range("I2").Value = 0

# Row 3
$ws3.Range("B3").Value = 0.1
$ws3.Range("C3").Value = 45
$ws3.Range("D3").Value = 66.7
$ws3.Range("F3").Value = 2
$ws3.Range("G3").Value = 2
$ws3.Range("H3").Value = 1
$ws3.Range("I3").Value = 0

# Row 4
$ws3.Range("B4").Value = 0.1
$ws3.Range("C4").Value = 85
$ws3.Range("D4").Value = 66.7
$ws3.Range("F4").Value = 2
$ws3.Range("G4").Value = 2
$ws3.Range("H4").Value = 1
$ws3.Range("I4").Value = 0

# Row 5
$ws3.Range("B5").Value = 0.1
$ws3.Range("C5").Value = 125
$ws3.Range("D5").Value = 66.7
$ws3.Range("F5").Value = 2
$ws3.Range("G5").Value = 2
$ws3.Range("H5").Value = 1
$ws3.Range("I5").Value = 0

# --- View/selection state for each sheet ---

# Sheet 1 "1. General": scroll so row 8 is the top-left row, select B24
[void]$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 8
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws1.Range("B24").Select()

# Sheet 2 "2. SPS - Area Specification": select D3, no longer the active tab
[void]$ws2.Activate()
[void]$ws2.Range("D3").Select()

# Sheet 3 "3. SPS - Generic Parameters": becomes the active/visible tab, select I6
[void]$ws3.Activate()
[void]$ws3.Range("I6").Select()
